$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.1
$ws.Range("H2").Value = 2.9
$ws.Range("I2").Value = 3.9
$ws.Range("L2").Value = 4.5
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.67
$ws.Range("W2").Value = 6
$ws.Range("X2").Value = 9
$ws.Range("Z2").Value = 19
$ws.Range("AN2").Value = 4
$ws.Range("AP2").Value = 29
$ws.Range("AU2").Value = 9
$ws.Range("AV2").Value = 67
$ws.Range("AY2").Value = 34
